# Updated symbol list on Sun Jan  8 18:49:15 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to the crypto table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'266.02"
$ws.Range("E2").Value = "'1.72%"
$ws.Range("D3").Value = "'26.74"
$ws.Range("E3").Value = "'-1.87%"
$ws.Range("D4").Value = "'4.696"
$ws.Range("E4").Value = "'-0.05%"
$ws.Range("D5").Value = "'0.06082"
$ws.Range("E5").Value = "'-1.75%"
$ws.Range("D6").Value = "'6.693"
$ws.Range("E6").Value = "'-0.37%"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("D8").Value = "'0.9051"
$ws.Range("E8").Value = "'-1.16%"
$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'0.02%"
$ws.Range("D10").Value = "'0.04933"
$ws.Range("E10").Value = "'5.39%"
$ws.Range("D11").Value = "'0.07106"
$ws.Range("E11").Value = "'0.39%"
$ws.Range("D12").Value = "'0.03129"
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("D13").Value = "'0.09019"
$ws.Range("E13").Value = "'-0.21%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'-0.61%"
$ws.Range("D15").Value = "'0.0006077"
$ws.Range("E15").Value = "'-1.43%"
$ws.Range("D16").Value = "'0.006159"
$ws.Range("E16").Value = "'1.47%"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("D18").Value = "'3.169"
$ws.Range("E18").Value = "'0.02%"
$ws.Range("D19").Value = "'2.278"
$ws.Range("E19").Value = "'3.79%"
$ws.Range("E20").Value = "'0.35%"
$ws.Range("E21").Value = "'-0.64%"
$ws.Range("D22").Value = "'4.101"
$ws.Range("E22").Value = "'-0.31%"
$ws.Range("D23").Value = "'0.04241"
$ws.Range("E23").Value = "'0.32%"
$ws.Range("D24").Value = "'0.001179"
$ws.Range("E24").Value = "'-3.02%"
$ws.Range("D25").Value = "'0.004134"
$ws.Range("E25").Value = "'8.72%"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E27").Value = "'5.05%"
$ws.Range("D40").Value = "'0.03918"
$ws.Range("E40").Value = "'-1.54%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.20%"
$ws.Range("D42").Value = "'0.004170"
$ws.Range("D43").Value = "'0.002112"
$ws.Range("E43").Value = "'-3.30%"
$ws.Range("E44").Value = "'-16.51%"
$ws.Range("D45").Value = "'0.00005129"
$ws.Range("E45").Value = "'-0.12%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D48").Value = "'0.1395"
$ws.Range("E48").Value = "'-16.28%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.07%"